# Auto-generated Excel COM-interop script to apply the diff
# "Update latest output (run 108)"

$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule": update summary schedule rows 2, 4, 5 ---
$schedule = $wb.Worksheets.Item("Schedule")
$schedule.Range("B2").Value = 46041.25
$schedule.Range("C2").Value = 6
$schedule.Range("D2").Value = 22.68
$schedule.Range("E2").Value = 567.7072635
$schedule.Range("F2").Value = 25.03118445767196
$schedule.Range("A4").Value = 46041.97916666666
$schedule.Range("B4").Value = 46042.14583333334
$schedule.Range("E4").Value = 449.6219714999999
$schedule.Range("F4").Value = 29.73690287698412
$schedule.Range("A5").Value = 46042.3125
$schedule.Range("C5").Value = 8.5
$schedule.Range("D5").Value = 32.13
$schedule.Range("E5").Value = -105.39780225
$schedule.Range("F5").Value = -3.280354878618114

# --- Sheet "Detailed": update price/status rows ---
$detailed = $wb.Worksheets.Item("Detailed")
$detailed.Range("E13").Value = "ON"
$detailed.Range("B37").Value = 0.45799
$detailed.Range("B38").Value = 11.89734
$detailed.Range("B39").Value = 0.01116
$detailed.Range("C39").Value = "historical"
$detailed.Range("B40").Value = 47.73914
$detailed.Range("C40").Value = "historical"
$detailed.Range("B41").Value = 57.3
$detailed.Range("B42").Value = 64.8901
$detailed.Range("B44").Value = 65
$detailed.Range("B45").Value = 65
$detailed.Range("B46").Value = 59.20766
$detailed.Range("B47").Value = 65
$detailed.Range("B48").Value = 59.77009
$detailed.Range("E48").Value = "OFF"
$detailed.Range("B49").Value = 61.73053
$detailed.Range("B50").Value = 57.06003
$detailed.Range("B51").Value = 57.06003
$detailed.Range("E56").Value = "ON"
$detailed.Range("B57").Value = 57.06003
$detailed.Range("B58").Value = 64.06528
$detailed.Range("B59").Value = 63.37713
$detailed.Range("B60").Value = 64.89
$detailed.Range("B61").Value = 68.88809999999999
$detailed.Range("B62").Value = 76.59223
$detailed.Range("B63").Value = 69.14248000000001
$detailed.Range("E64").Value = "OFF"
$detailed.Range("B65").Value = 4.78431
$detailed.Range("B66").Value = 0.05444
$detailed.Range("B67").Value = 0.0101
$detailed.Range("B68").Value = -2.54301
$detailed.Range("B70").Value = -7.17554
$detailed.Range("B71").Value = -7.11442
$detailed.Range("B72").Value = -7.93212
$detailed.Range("B73").Value = -7.92499
$detailed.Range("B74").Value = -6.91615
$detailed.Range("B75").Value = -9.99
$detailed.Range("B76").Value = -10
$detailed.Range("B77").Value = -11.01
$detailed.Range("B78").Value = -11.01
$detailed.Range("B79").Value = -11
$detailed.Range("B80").Value = -7.86626
$detailed.Range("B81").Value = -6.46667
$detailed.Range("B82").Value = -5.14727
$detailed.Range("B83").Value = -5.74039
$detailed.Range("B84").Value = -1.00044
$detailed.Range("B85").Value = 5.03982
$detailed.Range("B86").Value = 9.316129999999999
$detailed.Range("B87").Value = 33.11565
$detailed.Range("B88").Value = 56.98
$detailed.Range("B89").Value = 59.2087
$detailed.Range("B90").Value = 60.12181
$detailed.Range("B91").Value = 58.92354
$detailed.Range("B92").Value = 58.14193
$detailed.Range("B93").Value = 57.82439
$detailed.Range("B95").Value = 58.85642
$detailed.Range("B96").Value = 61.43326
